$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New N2 cell with a date/time number format (built-in format 22) ---
$ws.Range("N2").NumberFormat = "m/d/yy h:mm"

# --- D column: switch from _xlfn.CONCAT(...) to CONCATENATE(...) ---
$ws.Range("D4").Formula = '=PROPER(CONCATENATE(C4," ",B4))'
$ws.Range("D5:D38").Formula = '=PROPER(CONCATENATE(C5," ",B5))'

# --- G column: add YEARFRAC age-in-years formulas ---
$ws.Range("G4").Formula = "=YEARFRAC(F4,TODAY())"
$ws.Range("G5:G38").Formula = "=YEARFRAC(F5,TODAY())"

# --- I column: add "1 year later" date formulas (H + 365) ---
$ws.Range("I4").Formula = "=H4 + 365"
$ws.Range("I5:I38").Formula = "=H5 + 365"

# --- Update the active selection to the newly-filled I column ---
[void]$ws.Range("I4:I38").Select()
